$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the latest weekly report, pushing the existing
# historical rows (111-194) down to (113-196).
$ws.Rows("111:112").Insert()

# New row 111 - "Primera" quality record for the new reporting date.
$ws.Range("A111").Value = 9
$ws.Range("B111").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C111").Value = "Metropolitana"
$ws.Range("D111").Value2 = 45236
$ws.Range("E111").Value = 13
$ws.Range("F111").Value = 100114002
$ws.Range("G111").Value = "Camote"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 970
$ws.Range("K111").Value = 14000
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = 14485
$ws.Range("N111").Value = "$/caja 18 kilos"
$ws.Range("O111").Value = "Perú"
$ws.Range("P111").Value = 805
$ws.Range("Q111").Value = 18
$ws.Range("R111").Value = "Hortaliza"

# New row 112 - "Primera" quality record for the new reporting date.
$ws.Range("A112").Value = 9
$ws.Range("B112").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value2 = 45236
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 100114002
$ws.Range("G112").Value = "Camote"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 1060
$ws.Range("K112").Value = 12000
$ws.Range("L112").Value = 13000
$ws.Range("M112").Value = 12500
$ws.Range("N112").Value = "$/malla 18 kilos"
$ws.Range("O112").Value = "Perú"
$ws.Range("P112").Value = 694
$ws.Range("Q112").Value = 18
$ws.Range("R112").Value = "Hortaliza"
